$wb = $excel.ActiveWorkbook

# =========================================================================
# Numerical sheet: add "skewness" column (L) + refresh changed metric values
# =========================================================================
$numWs = $wb.Worksheets.Item("Numerical")

# Copy the existing bold/bordered header style (used by A1:K1) onto the new
# L1 header cell so it matches the rest of the header row.
$numWs.Range("K1").Copy()
$numWs.Range("L1").PasteSpecial(-4122)
$numWs.Range("L1").Value = "skewness"

$numericalUpdates = @(
  @{ row=2; vals=@{ B=13612; C=0; D=0; E=0; F=0; G=13612; H=1; I=15755; J=7841.89; K=7573.5; L=0.06 } },
  @{ row=3; vals=@{ B=5321; C=10; D=0; E=8291; F=60.91; G=2087; H=0; I=584625; J=20701.87; K=8500; L=5.71 } },
  @{ row=4; vals=@{ B=5323; C=5; D=0; E=8289; F=60.89; G=2731; H=0; I=2670982.5; J=34200.25; K=12000; L=12.29 } },
  @{ row=5; vals=@{ B=4938; C=1; D=0; E=8674; F=63.72; G=3491; H=0; I=984026.05; J=16807.41; K=5877.5; L=9.15 } },
  @{ row=6; vals=@{ B=4952; C=104; D=3; E=8660; F=63.62; G=1670; H=-5700; I=211833; J=2050.41; K=700; L=15.51 } },
  @{ row=7; vals=@{ B=7104; C=1781; D=0; E=6508; F=47.81; G=31; H=0; I=427; J=2.18; K=2; L=59.68 } },
  @{ row=8; vals=@{ B=7104; C=4148; D=0; E=6508; F=47.81; G=10; H=0; I=11; J=0.57; K=0; L=2.43 } },
  @{ row=9; vals=@{ B=7104; C=1781; D=0; E=6508; F=47.81; G=17; H=0; I=74; J=1.1; K=1; L=21.53 } },
  @{ row=10; vals=@{ B=7104; C=5762; D=0; E=6508; F=47.81; G=10; H=0; I=43; J=0.22; K=0; L=31.51 } },
  @{ row=11; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=12; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=13; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=14; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=15; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=16; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=17; vals=@{ B=0; C=0; D=0; E=13612; F=100; G=0 } },
  @{ row=18; vals=@{ B=396; C=0; D=0; E=13216; F=97.09; G=4; H=1; I=4; J=1.06; K=1; L=6.04 } },
  @{ row=19; vals=@{ B=396; C=257; D=0; E=13216; F=97.09; G=132; H=0; I=480000; J=11727.75; K=0; L=7.39 } },
  @{ row=20; vals=@{ B=1342; C=0; D=0; E=12270; F=90.14; G=13; H=1; I=235; J=1.68; K=1; L=35.23 } }
)

foreach ($item in $numericalUpdates) {
    $r = $item.row
    foreach ($col in $item.vals.Keys) {
        $numWs.Range("$col$r").Value = $item.vals[$col]
    }
}

# =========================================================================
# Categorical sheet: refresh changed counts/percentages
# =========================================================================
$catWs = $wb.Worksheets.Item("Categorical")

$categoricalUpdates = @(
  @{ row=2; vals=@{ B=5322; C=0; D=8290; E=60.9; F=5 } },
  @{ row=3; vals=@{ B=13598; C=0; D=14; E=0.1; F=8 } },
  @{ row=4; vals=@{ B=13598; C=0; D=14; E=0.1; F=6 } },
  @{ row=5; vals=@{ B=13576; C=0; D=36; E=0.26; F=5 } },
  @{ row=6; vals=@{ B=13332; C=0; D=280; E=2.06; F=4 } },
  @{ row=7; vals=@{ B=7104; C=0; D=6508; E=47.81; F=3 } },
  @{ row=8; vals=@{ B=13612; C=0; D=0; E=0; F=6 } },
  @{ row=9; vals=@{ B=13612; C=0; D=0; E=0; F=7 } },
  @{ row=10; vals=@{ B=13612; C=0; D=0; E=0; F=11 } }
)

foreach ($item in $categoricalUpdates) {
    $r = $item.row
    foreach ($col in $item.vals.Keys) {
        $catWs.Range("$col$r").Value = $item.vals[$col]
    }
}

# =========================================================================
# New "Usage" sheet (database usage / observability metadata)
# =========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$usageWs = $wb.Worksheets.Add($null, $lastSheet)
$usageWs.Name = "Usage"

$usageHeaders = @("schemaname","no_of_times_accessed","table_name","indexrelname","tables_usability","index_usability")

# Re-use the same bold/bordered header style as the other sheets' header rows.
$numWs.Range("A1").Copy()
$usageWs.Range("A1:F1").PasteSpecial(-4122)

for ($i = 0; $i -lt $usageHeaders.Length; $i++) {
    $usageWs.Cells.Item(1, $i + 1).Value = $usageHeaders[$i]
}

$usageWs.Cells.Item(2, 1).Value = "adaptiveai"
$usageWs.Cells.Item(2, 2).Value = 45
$usageWs.Cells.Item(2, 3).Value = "client_dim"
$usageWs.Cells.Item(2, 4).Value = ""
$usageWs.Cells.Item(2, 5).Value = "Used"
$usageWs.Cells.Item(2, 6).Value = "Index not used"
